$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set the new value in B2
$ws.Range("B2").Value = "S14523075"

# Clear the now-empty R2 cell (previously held only date formatting with no value)
$ws.Range("R2").Clear()

# Update the active selection to match the final state
$ws.Range("H19").Select()
